$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the original 'lock' column (column A). This shifts
# 'lock hex' -> A, 'fsu' -> B, 'start' -> C, 'end' -> D.
$ws.Columns.Item(1).Delete()

# Rename the (now) first column header from 'lock hex' to 'lock id'.
$ws.Range("A1").Value = "lock id"

# Strip the '0x' prefix from the hex lock id values so that, e.g.,
# '0xd716' becomes 'd716'.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null) {
        $text = [string]$val
        if ($text.StartsWith("0x")) {
            $cell.Value = $text.Substring(2)
        }
    }
}
